# Update the F column ("想去人数") values on the 展览, 演出 and 全部类型
# sheets from 0 to their real scraped values.

$wb = $excel.ActiveWorkbook

$values = @(356, 8, 10689, 331, 971, 140, 1318, 8253, 33, 464, 22, 216, 136, 3284, 39, 328, 769, 128, 1062, 286, 103, 1742)

# Sheet "展览" (exhibitions): rows 2..23
$ws1 = $wb.Worksheets.Item("展览")
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 6).Value = $values[$i]
}

# Sheet "演出" (performances): row 2 only
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 30

# Sheet "全部类型" (all types): rows 2..24 (same values as 展览, plus the 演出 value)
$ws4 = $wb.Worksheets.Item("全部类型")
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 6).Value = $values[$i]
}
$ws4.Cells.Item(24, 6).Value = 30
